$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new test-mail row (row 16) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(16, 1).Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(16, 3).Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$logs.Cells.Item(16, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item(16, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item(16, 6).Value = "2025-08-04 20:33:58"
$logs.Cells.Item(16, 7).Value = "Ja"
$logs.Cells.Item(16, 8).Value = "Ja"
$logs.Cells.Item(16, 9).Value = "Nee"
$logs.Cells.Item(16, 10).Value = "Nee"

# --- Extend the conditional formatting ranges to cover the new row ---
$ranges = @("D2:D16", "G2:G16", "H2:H16", "I2:I16", "J2:J16")
foreach ($sq in $ranges) {
    $target = $logs.Range($sq)
    $fc = $target.FormatConditions
    for ($i = 1; $i -le $fc.Count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($target)
    }
}

# --- Dashboard sheet: category counts re-sorted (Inkoop/Bestellingen now 3, overtaking Retour/Terugbetaling) ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(4, 2).Value = 3
$dash.Cells.Item(5, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(5, 2).Value = 2
